# Small correction on statistics slide:
# Slide 7 ("Some numbers") has a table; the second "Ratio calendar spam/spam"
# row label (row 4, col 1) should read "Ratio calendar spam/clean".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

# Locate the shape that holds the statistics table on this slide.
$tableShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTable) {
        $tableShape = $shape
        break
    }
}

$tbl = $tableShape.Table
$cell = $tbl.Cell(4, 1)
$cell.Shape.TextFrame.TextRange.Text = "Ratio calendar spam/clean"
